$wb = $excel.ActiveWorkbook

# Update "想去人数" (interest count) values on the 展览 sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 341
$ws1.Range("F4").Value = 1427

# Update the same cells on the 全部类型 sheet, which mirrors these rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 341
$ws4.Range("F4").Value = 1427
